# Apply updated crypto symbol/price data to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param(
        [string]$Address,
        [string]$Text
    )
    $rng = $ws.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value2 = $Text
}

# Price (column D) updates
Set-TextCell "D2"  "269.31"
Set-TextCell "D3"  "22.88"
Set-TextCell "D4"  "6.325"
Set-TextCell "D5"  "0.06197"
Set-TextCell "D6"  "3.641"
Set-TextCell "D7"  "6.689"
Set-TextCell "D8"  "1.393"
Set-TextCell "D9"  "0.8300"
Set-TextCell "D10" "0.01379"
Set-TextCell "D11" "0.1604"
Set-TextCell "D12" "0.08274"
Set-TextCell "D13" "0.03483"
Set-TextCell "D14" "0.03182"
Set-TextCell "D15" "0.09341"
Set-TextCell "D16" "3.837"
Set-TextCell "D17" "0.001660"
Set-TextCell "D19" "0.006418"
Set-TextCell "D20" "0.005671"
Set-TextCell "D23" "3.720"
Set-TextCell "D24" "2.324"
Set-TextCell "D27" "0.0002705"
Set-TextCell "D40" "0.04700"
Set-TextCell "D41" "0.006974"

# Rows 42 and 43 swap coin identity (BKEXToken <-> CEJI) with new prices/volume labels
Set-TextCell "B42" "CEJI"
Set-TextCell "C42" "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell "D42" "0.003801"
Set-TextCell "E42" "41CEJICEJIWorstin24h"

Set-TextCell "B43" "BKEXToken"
Set-TextCell "C43" "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell "D43" "0.1158"
Set-TextCell "E43" "42BKEXTokenBKK"

Set-TextCell "D44" "0.01191"
Set-TextCell "D45" "0.00006260"
Set-TextCell "D46" "0.0009901"
Set-TextCell "D48" "0.9202"
Set-TextCell "D49" "0.002195"
